$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'330.39"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'0.99%"
$c.Style = "Normal"

$c = $ws.Range("G2")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'44.29"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "'-0.61%"
$c.Style = "Normal"

$c = $ws.Range("G3")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.458"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'-2.44%"
$c.Style = "Normal"

$c = $ws.Range("G4")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.08025"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'-0.45%"
$c.Style = "Normal"

$c = $ws.Range("G5")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'1.992"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'4.83%"
$c.Style = "Normal"

$c = $ws.Range("G6")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.9535"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'0.90%"
$c.Style = "Normal"

$c = $ws.Range("G7")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'2.564"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'-4.62%"
$c.Style = "Normal"

$c = $ws.Range("G8")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.1141"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'-1.51%"
$c.Style = "Normal"

$c = $ws.Range("G9")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.1910"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'3.08%"
$c.Style = "Normal"

$c = $ws.Range("G10")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'10.57"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'25.16%"
$c.Style = "Normal"

$c = $ws.Range("G11")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.09893"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'-0.15%"
$c.Style = "Normal"

$c = $ws.Range("G12")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.04822"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'13.41%"
$c.Style = "Normal"

$c = $ws.Range("G13")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'-0.22%"
$c.Style = "Normal"

$c = $ws.Range("G14")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.001269"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'-0.14%"
$c.Style = "Normal"

$c = $ws.Range("G15")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'0.04073"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'-3.45%"
$c.Style = "Normal"

$c = $ws.Range("G16")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'0.005910"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'-0.35%"
$c.Style = "Normal"

$c = $ws.Range("G17")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B18")
$c.Value = "'LEO"
$c.Style = "Normal"

$c = $ws.Range("C18")
$c.Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'3.370"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'-6.64%"
$c.Style = "Normal"

$c = $ws.Range("G18")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B19")
$c.Value = "'GateToken"
$c.Style = "Normal"

$c = $ws.Range("C19")
$c.Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'4.396"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'2.12%"
$c.Style = "Normal"

$c = $ws.Range("G19")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B20")
$c.Value = "'BitpandaEcosystemToken"
$c.Style = "Normal"

$c = $ws.Range("C20")
$c.Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.3428"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'-1.99%"
$c.Style = "Normal"

$c = $ws.Range("G20")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B21")
$c.Value = "'ProBitToken"
$c.Style = "Normal"

$c = $ws.Range("C21")
$c.Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'0.1407"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'2.60%"
$c.Style = "Normal"

$c = $ws.Range("G21")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B22")
$c.Value = "'ZBToken"
$c.Style = "Normal"

$c = $ws.Range("C22")
$c.Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.2502"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'-0.23%"
$c.Style = "Normal"

$c = $ws.Range("G22")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B23")
$c.Value = "'BitKan"
$c.Style = "Normal"

$c = $ws.Range("C23")
$c.Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.001272"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'2.19%"
$c.Style = "Normal"

$c = $ws.Range("G23")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("B24")
$c.Value = "'HotbitToken"
$c.Style = "Normal"

$c = $ws.Range("C24")
$c.Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.004351"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'-2.53%"
$c.Style = "Normal"

$c = $ws.Range("G24")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'-5.02%"
$c.Style = "Normal"

$c = $ws.Range("G25")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.0003739"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'-6.38%"
$c.Style = "Normal"

$c = $ws.Range("G26")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G27")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G28")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G29")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G30")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G31")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G32")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G33")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G34")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G35")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G36")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("G37")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.02594"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'-1.70%"
$c.Style = "Normal"

$c = $ws.Range("G38")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.05795"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'6.11%"
$c.Style = "Normal"

$c = $ws.Range("G39")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.007539"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'-1.06%"
$c.Style = "Normal"

$c = $ws.Range("G40")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.1401"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'0.43%"
$c.Style = "Normal"

$c = $ws.Range("G41")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'0.03%"
$c.Style = "Normal"

$c = $ws.Range("G42")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'-1.86%"
$c.Style = "Normal"

$c = $ws.Range("G43")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.008838"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'-0.09%"
$c.Style = "Normal"

$c = $ws.Range("G44")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.00007119"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'2.75%"
$c.Style = "Normal"

$c = $ws.Range("G45")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'-0.27%"
$c.Style = "Normal"

$c = $ws.Range("G46")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.0005792"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'-0.33%"
$c.Style = "Normal"

$c = $ws.Range("G47")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.003525"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'55.11%"
$c.Style = "Normal"

$c = $ws.Range("G48")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'0.003494"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'-4.70%"
$c.Style = "Normal"

$c = $ws.Range("G49")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'-0.27%"
$c.Style = "Normal"

$c = $ws.Range("G50")
$c.Value = "'10"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'-0.27%"
$c.Style = "Normal"

$c = $ws.Range("G51")
$c.Value = "'10"
$c.Style = "Normal"
